# Auto-generated edit script: refresh computed market-price / profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all 8 job sheets,
# matching a scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2496.1667
$ws.Range("I18").Value = 3245
$ws.Range("K18").Value = 3245
$ws.Range("M18").Value = -2961
$ws.Range("H41").Value = 542
$ws.Range("J41").Value = 1822.6666
$ws.Range("L41").Value = 1822.6666
$ws.Range("N41").Value = -2702.6666
$ws.Range("H62").Value = 20333
$ws.Range("J62").Value = 6000
$ws.Range("L62").Value = 6000
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 20333
$ws.Range("J65").Value = 6000
$ws.Range("L65").Value = 30000
$ws.Range("N65").Value = -36240
$ws.Range("H80").Value = 4643.6924
$ws.Range("I80").Value = 3827.7144
$ws.Range("J80").Value = 5595.6665
$ws.Range("K80").Value = 11483.1432
$ws.Range("L80").Value = 16786.9995
$ws.Range("M80").Value = -10485.1432
$ws.Range("N80").Value = -18782.9995
$ws.Range("H83").Value = 4643.6924
$ws.Range("I83").Value = 3827.7144
$ws.Range("J83").Value = 5595.6665
$ws.Range("K83").Value = 34449.4296
$ws.Range("L83").Value = 50360.9985
$ws.Range("M83").Value = -29457.4296
$ws.Range("N83").Value = -60344.9985
$ws.Range("H116").Value = 22830.8
$ws.Range("I116").Value = 6285.375
$ws.Range("K116").Value = 6285.375
$ws.Range("M116").Value = -2843.375
$ws.Range("H138").Value = 3912.1272
$ws.Range("J138").Value = 3906.16
$ws.Range("L138").Value = 11718.48
$ws.Range("N138").Value = -21998.48

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1050.1428
$ws.Range("I2").Value = 1114.9667
$ws.Range("J2").Value = 661.2
$ws.Range("K2").Value = 1114.9667
$ws.Range("L2").Value = 661.2
$ws.Range("M2").Value = -1001.9667
$ws.Range("N2").Value = -887.2
$ws.Range("H32").Value = 213379.89
$ws.Range("I32").Value = 267911.03
$ws.Range("J32").Value = 24999.637
$ws.Range("K32").Value = 267911.03
$ws.Range("L32").Value = 24999.637
$ws.Range("M32").Value = -267624.03
$ws.Range("N32").Value = -25573.637
$ws.Range("H37").Value = 1263
$ws.Range("J37").Value = 1263
$ws.Range("L37").Value = 1263
$ws.Range("N37").Value = -1809
$ws.Range("H44").Value = 34000
$ws.Range("J44").Value = 34000
$ws.Range("L44").Value = 34000
$ws.Range("N44").Value = -34976
$ws.Range("H45").Value = 1744.1111
$ws.Range("I45").Value = 1539.4
$ws.Range("K45").Value = 1539.4
$ws.Range("M45").Value = -1162.4
$ws.Range("H55").Value = 35000
$ws.Range("J55").Value = 35000
$ws.Range("L55").Value = 35000
$ws.Range("N55").Value = -35630
$ws.Range("H63").Value = 9469.833000000001
$ws.Range("I63").Value = 9205
$ws.Range("K63").Value = 9205
$ws.Range("M63").Value = -8519
$ws.Range("H66").Value = 9469.833000000001
$ws.Range("I66").Value = 9205
$ws.Range("K66").Value = 46025
$ws.Range("M66").Value = -42593
$ws.Range("H94").Value = 48990
$ws.Range("J94").Value = 48990
$ws.Range("L94").Value = 48990
$ws.Range("N94").Value = -50792
$ws.Range("H116").Value = 1050.1428
$ws.Range("I116").Value = 1114.9667
$ws.Range("J116").Value = 661.2
$ws.Range("K116").Value = 1114.9667
$ws.Range("L116").Value = 661.2
$ws.Range("M116").Value = 1179.0333
$ws.Range("N116").Value = -5249.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1050.1428
$ws.Range("I3").Value = 1114.9667
$ws.Range("J3").Value = 661.2
$ws.Range("K3").Value = 1114.9667
$ws.Range("L3").Value = 661.2
$ws.Range("M3").Value = -1000.9667
$ws.Range("N3").Value = -889.2
$ws.Range("H82").Value = 18509.934
$ws.Range("H85").Value = 18509.934
$ws.Range("H99").Value = 13030.923
$ws.Range("I99").Value = 16415.9
$ws.Range("K99").Value = 16415.9
$ws.Range("M99").Value = -14917.9

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 23866.238
$ws.Range("I7").Value = 100045.2
$ws.Range("K7").Value = 100045.2
$ws.Range("M7").Value = -99932.2
$ws.Range("H17").Value = 25250
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 50000
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 50000
$ws.Range("M17").Value = -326
$ws.Range("N17").Value = -50348
$ws.Range("H22").Value = 46820.363
$ws.Range("J22").Value = 72963.14
$ws.Range("L22").Value = 72963.14
$ws.Range("N22").Value = -73663.14
$ws.Range("H31").Value = 2781449.2
$ws.Range("I31").Value = 4275550.5
$ws.Range("K31").Value = 4275550.5
$ws.Range("M31").Value = -4275255.5
$ws.Range("H34").Value = 2781449.2
$ws.Range("I34").Value = 4275550.5
$ws.Range("K34").Value = 4275550.5
$ws.Range("M34").Value = -4275348.5
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value = 19507.445
$ws.Range("I93").Value = 19507.445
$ws.Range("K93").Value = 19507.445
$ws.Range("M93").Value = -17635.445
$ws.Range("H103").Value = 14205.154
$ws.Range("I103").Value = 14205.154
$ws.Range("K103").Value = 14205.154
$ws.Range("M103").Value = -13033.154
$ws.Range("H122").Value = 6653
$ws.Range("I122").Value = 1738.0646
$ws.Range("K122").Value = 5214.1938
$ws.Range("M122").Value = -2764.1938

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 188
$ws.Range("J2").Value = 350.33334
$ws.Range("L2").Value = 2102.00004
$ws.Range("N2").Value = -2328.00004
$ws.Range("H105").Value = 25000
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 75000
$ws.Range("N105").Value = -80242
$ws.Range("H113").Value = 2500.9375
$ws.Range("I113").Value = 1291.5
$ws.Range("J113").Value = 2673.7144
$ws.Range("K113").Value = 3874.5
$ws.Range("L113").Value = 8021.1432
$ws.Range("M113").Value = -1704.5
$ws.Range("N113").Value = -12361.1432
$ws.Range("H132").Value = 4662.4062
$ws.Range("J132").Value = 5266.778
$ws.Range("L132").Value = 47401.002
$ws.Range("N132").Value = -52461.002
$ws.Range("H141").Value = 2235.5
$ws.Range("I141").Value = 2235.5
$ws.Range("K141").Value = 6706.5
$ws.Range("M141").Value = -1526.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 5126.5835
$ws.Range("I21").Value = 5156.364
$ws.Range("K21").Value = 5156.364
$ws.Range("M21").Value = -4983.364
$ws.Range("H30").Value = 5126.5835
$ws.Range("I30").Value = 5156.364
$ws.Range("K30").Value = 5156.364
$ws.Range("M30").Value = -5051.364
$ws.Range("H62").Value = 44946
$ws.Range("I62").Value = 44946
$ws.Range("K62").Value = 44946
$ws.Range("M62").Value = -44260
$ws.Range("H65").Value = 44946
$ws.Range("I65").Value = 44946
$ws.Range("K65").Value = 134838
$ws.Range("M65").Value = -131406
$ws.Range("H113").Value = 3052.25
$ws.Range("I113").Value = 3052.25
$ws.Range("K113").Value = 3052.25
$ws.Range("M113").Value = -882.25
$ws.Range("H126").Value = 14003.333
$ws.Range("I126").Value = 17822.727
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 53468.181
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -50998.181
$ws.Range("N126").Value = -15440

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 4964.25
$ws.Range("J13").Value = 4941.5
$ws.Range("L13").Value = 4941.5
$ws.Range("N13").Value = -5221.5
$ws.Range("H23").Value = 9680.154
$ws.Range("I23").Value = 9677.727999999999
$ws.Range("K23").Value = 9677.727999999999
$ws.Range("M23").Value = -9447.727999999999
$ws.Range("H69").Value = 39000
$ws.Range("J69").Value = 39000
$ws.Range("L69").Value = 39000
$ws.Range("N69").Value = -40622
$ws.Range("H72").Value = 39000
$ws.Range("J72").Value = 39000
$ws.Range("L72").Value = 117000
$ws.Range("N72").Value = -125112
$ws.Range("H132").Value = 3772340
$ws.Range("I132").Value = 6876075.5
$ws.Range("J132").Value = 3518.0715
$ws.Range("K132").Value = 20628226.5
$ws.Range("L132").Value = 10554.2145
$ws.Range("M132").Value = -20625696.5
$ws.Range("N132").Value = -15614.2145
$ws.Range("H136").Value = 8629339
$ws.Range("I136").Value = 15631663
$ws.Range("J136").Value = 5961787
$ws.Range("K136").Value = 46894989
$ws.Range("L136").Value = 17885361
$ws.Range("M136").Value = -46892439
$ws.Range("N136").Value = -17890461

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 84989.78999999999
$ws.Range("I122").Value = 6129.3
$ws.Range("K122").Value = 18387.9
$ws.Range("M122").Value = -15937.9
$ws.Range("H132").Value = 4763914.5
$ws.Range("J132").Value = 2575
$ws.Range("L132").Value = 7725
$ws.Range("N132").Value = -12785
$ws.Range("H136").Value = 34783280
$ws.Range("I136").Value = 7247160.5
$ws.Range("J136").Value = 200000000
$ws.Range("K136").Value = 21741481.5
$ws.Range("L136").Value = 600000000
$ws.Range("M136").Value = -21738931.5
$ws.Range("N136").Value = -600005100

